# Weekly fruit/vegetable price update.
# Inserts 4 new observation rows at the top of the "Poroto granado" block
# (pushing the existing rows 408-426 down to 412-430) and populates the
# new rows 408-411 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 408:426 down to 412:430, leaving 408:411 blank
# (Excel copies the row-408 formatting into the newly inserted rows,
# which keeps column D's date number format intact).
$ws.Rows("408:411").Insert()

# --- Row 408 --------------------------------------------------------
$ws.Range("A408").Value = 6
$ws.Range("B408").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C408").Value = "Metropolitana"
$ws.Range("D408").Value = 44610
$ws.Range("E408").Value = 13
$ws.Range("F408").Value = 100112030
$ws.Range("G408").Value = "Poroto granado"
$ws.Range("H408").Value = "Sin especificar"
$ws.Range("I408").Value = "Primera"
$ws.Range("J408").Value = 1400
$ws.Range("K408").Value = 20000
$ws.Range("L408").Value = 23000
$ws.Range("M408").Value = 21286
$ws.Range("N408").Value = "$/saco 25 kilos"
$ws.Range("O408").Value = "Región Metropolitana"
$ws.Range("P408").Value = 851
$ws.Range("Q408").Value = 25
$ws.Range("R408").Value = "Hortaliza"

# --- Row 409 --------------------------------------------------------
$ws.Range("A409").Value = 6
$ws.Range("B409").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C409").Value = "Metropolitana"
$ws.Range("D409").Value = 44610
$ws.Range("E409").Value = 13
$ws.Range("F409").Value = 100112030
$ws.Range("G409").Value = "Poroto granado"
$ws.Range("H409").Value = "Sin especificar"
$ws.Range("I409").Value = "Primera"
$ws.Range("J409").Value = 900
$ws.Range("K409").Value = 18000
$ws.Range("L409").Value = 20000
$ws.Range("M409").Value = 19111
$ws.Range("N409").Value = "$/saco 25 kilos"
$ws.Range("O409").Value = "Región del Maule"
$ws.Range("P409").Value = 764
$ws.Range("Q409").Value = 25
$ws.Range("R409").Value = "Hortaliza"

# --- Row 410 --------------------------------------------------------
$ws.Range("A410").Value = 6
$ws.Range("B410").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C410").Value = "Metropolitana"
$ws.Range("D410").Value = 44610
$ws.Range("E410").Value = 13
$ws.Range("F410").Value = 100112030
$ws.Range("G410").Value = "Poroto granado"
$ws.Range("H410").Value = "Sin especificar"
$ws.Range("I410").Value = "Segunda"
$ws.Range("J410").Value = 400
$ws.Range("K410").Value = 17000
$ws.Range("L410").Value = 17000
$ws.Range("M410").Value = 17000
$ws.Range("N410").Value = "$/saco 25 kilos"
$ws.Range("O410").Value = "Región Metropolitana"
$ws.Range("P410").Value = 680
$ws.Range("Q410").Value = 25
$ws.Range("R410").Value = "Hortaliza"

# --- Row 411 --------------------------------------------------------
$ws.Range("A411").Value = 6
$ws.Range("B411").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C411").Value = "Metropolitana"
$ws.Range("D411").Value = 44610
$ws.Range("E411").Value = 13
$ws.Range("F411").Value = 100112030
$ws.Range("G411").Value = "Poroto granado"
$ws.Range("H411").Value = "Sin especificar"
$ws.Range("I411").Value = "Segunda"
$ws.Range("J411").Value = 200
$ws.Range("K411").Value = 15000
$ws.Range("L411").Value = 15000
$ws.Range("M411").Value = 15000
$ws.Range("N411").Value = "$/saco 25 kilos"
$ws.Range("O411").Value = "Región del Maule"
$ws.Range("P411").Value = 600
$ws.Range("Q411").Value = 25
$ws.Range("R411").Value = "Hortaliza"
